$wb = $excel.ActiveWorkbook

function Set-TextValue($cell, $text) {
    # Force a numeric-looking string to be stored as genuine text (matches
    # the source data, which keeps things like fund codes / percentages as
    # text rather than numbers), then strip the "quote prefix" styling that
    # Excel normally stamps on such cells so the cell format stays plain.
    $cell.NumberFormat = "@"
    $cell.Value2 = $text
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------------
# 1. Update the "总计" (summary) sheet: insert a new "2022-Q4" row at the top
#    of the data (row 2), pushing the existing rows down by one.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Existing data rows (2..8) hold quarters 2022-Q3 .. 2020-Q4. They all need
# to move down one row (new rows 3..9). Walk bottom-up so we never clobber
# a value before reading it. (Use Value2 -- Value's getter is unreliable in
# this host.)
for ($r = 8; $r -ge 2; $r--) {
    $dst = $r + 1
    $summary.Cells.Item($dst, 2).Value2 = $summary.Cells.Item($r, 2).Value2
    $summary.Cells.Item($dst, 3).Value2 = $summary.Cells.Item($r, 3).Value2
    $summary.Cells.Item($dst, 4).Value2 = $summary.Cells.Item($r, 4).Value2
}

# Row 2 becomes the new 2022-Q4 entry.
$summary.Cells.Item(2, 2).Value2 = "2022-Q4"
$summary.Cells.Item(2, 3).Value2 = 20
$summary.Cells.Item(2, 4).Value2 = 5.82

# Column A is just a running 0-based index; refresh it for all 8 data rows.
for ($r = 2; $r -le 9; $r++) {
    $summary.Cells.Item($r, 1).Value2 = $r - 2
}
$summary.Range("A2").Copy()
$summary.Range("A3:A9").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Insert a brand-new "2022-Q4" sheet (fund holdings detail) right after
#    "总计", matching the position/order used elsewhere in the workbook.
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $summary)
$q4.Name = "2022-Q4"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 2; $c -le 8; $c++) {
    Set-TextValue $q4.Cells.Item(1, $c) $headers[$c - 2]
}

$rows = @(
    @(0, "014591", "广发瑞誉一年持有期混合A", "35.89", "92.43", "4.22", "1.5146", 7),
    @(1, "513090", "易方达中证香港证券投资主题ETF", "11.28", "97.07", "12.47", "1.4066", 3),
    @(2, "014362", "睿远稳进配置两年持有混合A", "64.57", "37.12", "1.25", "0.8071", 9),
    @(3, "014363", "睿远稳进配置两年持有混合C", "35.61", "37.12", "1.25", "0.4451", 9),
    @(4, "012943", "广发稳睿六个月持有期混合A", "20.00", "26.98", "2.01", "0.4020", 6),
    @(5, "012944", "广发稳睿六个月持有期混合C", "18.37", "26.98", "2.01", "0.3692", 6),
    @(6, "202801", "南方全球精选配置（QDII-FOF）", "17.02", "32.64", "1.56", "0.2655", 6),
    @(7, "014592", "广发瑞誉一年持有期混合C", "4.26", "92.43", "4.22", "0.1798", 7),
    @(8, "007109", "南方沪港深核心优势混合", "1.99", "89.23", "4.35", "0.0866", 7),
    @(9, "003413", "华泰柏瑞新经济沪港深混合", "1.44", "94.26", "5.10", "0.0734", 9),
    @(10, "013659", "中融金融鑫选3个月持有混合A", "1.39", "90.62", "4.14", "0.0575", 10),
    @(11, "011355", "华泰柏瑞港股通时代机遇混合A", "0.70", "94.61", "6.09", "0.0426", 8),
    @(12, "013660", "中融金融鑫选3个月持有混合C", "0.83", "90.62", "4.14", "0.0344", 10),
    @(13, "011969", "建信港股通精选混合A", "0.62", "87.38", "4.99", "0.0309", 6),
    @(14, "005576", "华泰柏瑞新金融地产灵活配置混合A", "0.77", "94.17", "3.86", "0.0297", 9),
    @(15, "460010", "华泰柏瑞亚洲领导企业混合（QDII）", "0.52", "97.17", "5.67", "0.0295", 6),
    @(16, "011356", "华泰柏瑞港股通时代机遇混合C", "0.39", "94.61", "6.09", "0.0238", 8),
    @(17, "011970", "建信港股通精选混合C", "0.24", "87.38", "4.99", "0.0120", 6),
    @(18, "016374", "华泰柏瑞新金融地产灵活配置混合C", "0.22", "94.17", "3.86", "0.0085", 9),
    @(19, "003279", "融通沪港深智慧生活灵活配置混合", "0.05", "53.96", "4.98", "0.0025", 5)
)

$r = 2
foreach ($row in $rows) {
    $q4.Cells.Item($r, 1).Value2 = $row[0]
    Set-TextValue $q4.Cells.Item($r, 2) $row[1]
    Set-TextValue $q4.Cells.Item($r, 3) $row[2]
    Set-TextValue $q4.Cells.Item($r, 4) $row[3]
    Set-TextValue $q4.Cells.Item($r, 5) $row[4]
    Set-TextValue $q4.Cells.Item($r, 6) $row[5]
    Set-TextValue $q4.Cells.Item($r, 7) $row[6]
    $q4.Cells.Item($r, 8).Value2 = $row[7]
    $r++
}

# Match the bold/centred/bordered style ("s=2") used for the header row and
# the index column (A) on every other sheet in this workbook.
$summary.Range("B1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$summary.Range("A2").Copy()
$q4.Range("A2:A21").PasteSpecial(-4122)
